# Orientações para slides - edits:
#  1) "- " + "Para o arquivo Índice Global de Inovação 2019:" paragraph ->
#     merged into a single bold run "- Para o arquivo Índice Global de
#     Inovação 2019:", paragraph gets a left indent (708 twips) and bold
#     paragraph-mark formatting.
#  2) "7º slide – Principal cluster..." paragraph -> all runs highlighted
#     yellow.

$d = $word.ActiveDocument

$xmlHeader = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>'
$xmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- Locate the two target paragraphs by their distinctive text ---
$paraIndice = $null
$para7Slide = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*Para o arquivo*ndice Global de Inova*o 2019:*") {
        $paraIndice = $p
    }
    if ($t -like "*Principal cluster de economias*") {
        $para7Slide = $p
    }
}

# --- Edit 1: "Para o arquivo Índice Global de Inovação 2019:" paragraph ---
$r1 = $paraIndice.Range
$frag1 = $xmlHeader + `
    '<w:p w14:paraId="274658F5" w14:textId="31AC4DCF" w:rsidR="003F212E" w:rsidRDefault="003F212E">' + `
        '<w:pPr>' + `
            '<w:ind w:left="708"/>' + `
            '<w:rPr><w:b/><w:bCs/></w:rPr>' + `
        '</w:pPr>' + `
        '<w:r w:rsidRPr="00B32E65">' + `
            '<w:rPr><w:b/><w:bCs/></w:rPr>' + `
            '<w:t>- Para o arquivo Índice Global de Inovação 2019:</w:t>' + `
        '</w:r>' + `
    '</w:p>' + $xmlFooter
$r1.InsertXML($frag1)

# --- Edit 2: "7º slide – Principal cluster..." paragraph gets yellow highlight ---
$r2 = $para7Slide.Range
$frag2 = $xmlHeader + `
    '<w:p w14:paraId="31921434" w14:textId="23367336" w:rsidR="00B32E65" w:rsidRDefault="00B32E65">' + `
        '<w:r>' + `
            '<w:rPr><w:highlight w:val="yellow"/></w:rPr>' + `
            '<w:t xml:space="preserve">7º slide – Principal cluster de economias ou...entre os 50 mais importantes, 2019 (tabela C – </w:t>' + `
        '</w:r>' + `
        '<w:proofErr w:type="spellStart"/>' + `
        '<w:r>' + `
            '<w:rPr><w:highlight w:val="yellow"/></w:rPr>' + `
            '<w:t>pg</w:t>' + `
        '</w:r>' + `
        '<w:proofErr w:type="spellEnd"/>' + `
        '<w:r>' + `
            '<w:rPr><w:highlight w:val="yellow"/></w:rPr>' + `
            '<w:t xml:space="preserve"> 21)</w:t>' + `
        '</w:r>' + `
    '</w:p>' + $xmlFooter
$r2.InsertXML($frag2)

Write-Output "Applied both edits."
